# Add two new columns ("near school" and "ratings") before the existing
# "dorm" column, shifting dorm/rating/location/barangay/distance/rides/details
# two columns to the right (G:H -> I:J, etc.), then update the header labels
# and data values to match the new dataset upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at G:H; everything from G onward shifts right by 2.
$ws.Range("G1:H1").EntireColumn.Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 7).Value = "near school"   # G1
$ws.Cells.Item(1, 8).Value = "ratings"       # H1
# I1:O1 already hold the shifted former values (dorm, rating, location,
# barangay, distance, rides, details) thanks to the column insert.

# --- Row 2 (Kim Legendary Living) ---
$ws.Cells.Item(2, 1).Value = 12345                      # A2 budget
$ws.Cells.Item(2, 6).Value = "Upper Bicutan"             # F2 barangay
$ws.Cells.Item(2, 7).Value = "N/A"                       # G2 near school
$ws.Cells.Item(2, 8).Value = "{}"                        # H2 ratings
$ws.Cells.Item(2, 9).Value = "Kim Legendary Living"      # I2 dorm
$ws.Cells.Item(2, 10).Value = 1                          # J2 rating
$ws.Cells.Item(2, 11).Value = "Taguig"                   # K2 location
$ws.Cells.Item(2, 12).Value = "Upper Bicutan"            # L2 barangay
$ws.Cells.Item(2, 13).Value = 1.5                        # M2 distance
$ws.Cells.Item(2, 14).Value = 2                          # N2 rides
$ws.Cells.Item(2, 15).Value = "{'name': 'Kim Legendary Living', 'budget': 12345, 'location': 'Taguig', 'room_size': 'Single', 'bathroom': 'Private', 'environment': 'Social', 'distance': 1.5, 'rides': 2, 'barangay': 'Upper Bicutan', 'address': 'Blk 69 Lot 69 Tibi St', 'landmark': 'JR Store', 'near_school': 'TCU', 'link': 'https://web.facebook.com/groups/982859945207269/posts/2590425444450703/', 'wifi': 'Available', 'likes': nan}"

# --- Row 3 (Serendipity Living) ---
$ws.Cells.Item(3, 7).Value = ""                          # G3 near school (blank)
$ws.Cells.Item(3, 8).Value = ""                          # H3 ratings (blank)
$ws.Cells.Item(3, 9).Value = "Serendipity Living"        # I3 dorm
$ws.Cells.Item(3, 10).Value = 0.000000007908262509965794 # J3 rating
$ws.Cells.Item(3, 11).Value = "Paranaque"                # K3 location
$ws.Cells.Item(3, 12).Value = "Hagonoy"                  # L3 barangay
$ws.Cells.Item(3, 13).Value = 2.3                        # M3 distance
$ws.Cells.Item(3, 14).Value = 3                          # N3 rides
$ws.Cells.Item(3, 15).Value = "{'name': 'Serendipity Living', 'budget': 1100, 'location': 'Paranaque', 'room_size': 'Double', 'bathroom': 'Private', 'environment': 'Quiet', 'distance': 2.3, 'rides': 3, 'barangay': 'Hagonoy', 'address': nan, 'landmark': nan, 'near_school': 'Harmony High', 'link': nan, 'wifi': 'Not Available', 'likes': nan}"

# --- Row 4 (Pinnacle Residence) ---
$ws.Cells.Item(4, 7).Value = ""                          # G4 near school (blank)
$ws.Cells.Item(4, 8).Value = ""                          # H4 ratings (blank)
$ws.Cells.Item(4, 9).Value = "Pinnacle Residence"        # I4 dorm
$ws.Cells.Item(4, 10).Value = 0.000000007908262509965794 # J4 rating
$ws.Cells.Item(4, 11).Value = "Paranaque"                # K4 location
$ws.Cells.Item(4, 12).Value = "Lower Bicutan"            # L4 barangay
$ws.Cells.Item(4, 13).Value = 2                          # M4 distance
$ws.Cells.Item(4, 14).Value = 2                          # N4 rides
$ws.Cells.Item(4, 15).Value = "{'name': 'Pinnacle Residence', 'budget': 1100, 'location': 'Paranaque', 'room_size': 'Double', 'bathroom': 'Private', 'environment': 'Quiet', 'distance': 2.0, 'rides': 2, 'barangay': 'Lower Bicutan', 'address': nan, 'landmark': nan, 'near_school': 'City Heights Academy', 'link': nan, 'wifi': 'Not Available', 'likes': nan}"
